$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Looks-Numeric {
    param([string]$value)
    return $value -match '^[+-]?[0-9]*\.?[0-9]+$'
}

function Set-TextValue {
    param($cell, [string]$value)
    if (Looks-Numeric $value) {
        # Force Excel to store the value as text (not re-interpret as a
        # number) by using the classic leading-apostrophe text prefix.
        $cell.Formula = "'" + $value
    } else {
        $cell.Value = $value
    }
}

Set-TextValue $ws.Range("D2") '52.305.12'
Set-TextValue $ws.Range("E2") '  +1.22%  '
Set-TextValue $ws.Range("D3") '2.838.38'
Set-TextValue $ws.Range("E3") '  +3.36%  '
Set-TextValue $ws.Range("D4") '1.00'
Set-TextValue $ws.Range("E4") '  +0.03%  '
Set-TextValue $ws.Range("D5") '357.29'
Set-TextValue $ws.Range("E5") '  +7.43%  '
Set-TextValue $ws.Range("D6") '114.76'
Set-TextValue $ws.Range("E6") '  -1.63%  '
Set-TextValue $ws.Range("D7") '0.548'
Set-TextValue $ws.Range("E7") '  +2.72%  '
Set-TextValue $ws.Range("E8") '  +0.06%  '
Set-TextValue $ws.Range("D9") '0.607'
Set-TextValue $ws.Range("E9") '  +5.92%  '
Set-TextValue $ws.Range("D10") '42.11'
Set-TextValue $ws.Range("E10") '  +1.39%  '
Set-TextValue $ws.Range("E11") '  +1.90%  '
Set-TextValue $ws.Range("D12") '20.06'
Set-TextValue $ws.Range("E12") '  -0.24%  '
Set-TextValue $ws.Range("E13") '  +1.33%  '
Set-TextValue $ws.Range("D14") '7.80'
Set-TextValue $ws.Range("E14") '  +3.15%  '
Set-TextValue $ws.Range("D15") '3.282.45'
Set-TextValue $ws.Range("E15") '  +3.46%  '
Set-TextValue $ws.Range("D16") '2.827.56'
Set-TextValue $ws.Range("E16") '  +3.06%  '
Set-TextValue $ws.Range("D17") '0.894'
Set-TextValue $ws.Range("E17") '  +1.54%  '
Set-TextValue $ws.Range("D18") '52.240.74'
Set-TextValue $ws.Range("E18") '  +1.29%  '
Set-TextValue $ws.Range("E19") '  +2.11%  '
Set-TextValue $ws.Range("E20") '  +7.07%  '
Set-TextValue $ws.Range("E21") '  +1.97%  '
Set-TextValue $ws.Range("E22") '  +3.00%  '
Set-TextValue $ws.Range("D23") '271.04'
Set-TextValue $ws.Range("E23") '  -3.02%  '
Set-TextValue $ws.Range("D24") '69.77'
Set-TextValue $ws.Range("E24") '  +0.34%  '
Set-TextValue $ws.Range("E25") '  +6.35%  '
Set-TextValue $ws.Range("D26") '26.86'
Set-TextValue $ws.Range("E26") '  +0.52%  '
Set-TextValue $ws.Range("D28") '10.26'
Set-TextValue $ws.Range("E28") '  +0.74%  '
Set-TextValue $ws.Range("E29") '  +1.32%  '
Set-TextValue $ws.Range("E30") '  +0.92%  '
Set-TextValue $ws.Range("D31") '50.68'
Set-TextValue $ws.Range("E31") '  +0.97%  '
Set-TextValue $ws.Range("D32") '33.92'
Set-TextValue $ws.Range("E32") '  -3.16%  '
Set-TextValue $ws.Range("D33") '5.86'
Set-TextValue $ws.Range("E33") '  +5.62%  '
Set-TextValue $ws.Range("D34") '0.0439'
Set-TextValue $ws.Range("E34") '  +27.03%  '
Set-TextValue $ws.Range("D35") '0.0831'
Set-TextValue $ws.Range("E35") '  +1.39%  '
Set-TextValue $ws.Range("E36") '  -0.02%  '
Set-TextValue $ws.Range("E37") '  +1.17%  '
Set-TextValue $ws.Range("E38") '  -0.55%  '
Set-TextValue $ws.Range("D39") '18.51'
Set-TextValue $ws.Range("E39") '  -2.78%  '
Set-TextValue $ws.Range("E40") '  +2.07%  '
Set-TextValue $ws.Range("B41") 'EnergySwap'
Set-TextValue $ws.Range("C41") 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range("D41") '23.72'
Set-TextValue $ws.Range("E41") '  +2.38%  '
Set-TextValue $ws.Range("B42") 'Stacks'
Set-TextValue $ws.Range("C42") 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range("D42") '2.56'
Set-TextValue $ws.Range("E42") '  +7.70%  '
Set-TextValue $ws.Range("D43") '0.115'
Set-TextValue $ws.Range("E43") '  +2.00%  '
Set-TextValue $ws.Range("D44") '126.45'
Set-TextValue $ws.Range("E44") '  -1.83%  '
Set-TextValue $ws.Range("E45") '  +2.21%  '
Set-TextValue $ws.Range("E46") '  +1.96%  '
Set-TextValue $ws.Range("D47") '2.047.87'
Set-TextValue $ws.Range("E47") '  -2.40%  '
Set-TextValue $ws.Range("E48") '  +3.80%  '
Set-TextValue $ws.Range("B49") 'SEI'
Set-TextValue $ws.Range("C49") 'https://coinranking.com/coin/8nxCqs-uj+sei-sei'
Set-TextValue $ws.Range("D49") '0.966'
Set-TextValue $ws.Range("E49") '  +11.93%  '
Set-TextValue $ws.Range("B50") 'THORChain'
Set-TextValue $ws.Range("C50") 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue $ws.Range("D50") '5.74'
Set-TextValue $ws.Range("E50") '  +4.01%  '
Set-TextValue $ws.Range("D51") '60.30'
Set-TextValue $ws.Range("E51") '  +0.71%  '
